$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.440.70'
$ws.Range('E2').Value = '  +0.29%  '
$ws.Range('D3').Value = '1.940.18'
$ws.Range('E3').Value = '  +0.23%  '
$ws.Range('D4').Value = "'1.007"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.62%  '
$ws.Range('D5').Value = "'0.7524"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +3.82%  '
$ws.Range('D6').Value = "'246.08"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.42%  '
$ws.Range('E7').Value = '  +0.54%  '
$ws.Range('D8').Value = "'0.3187"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -3.50%  '
$ws.Range('D9').Value = "'27.66"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.96%  '
$ws.Range('D10').Value = "'0.06994"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -3.40%  '
$ws.Range('D11').Value = "'0.7808"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -3.36%  '
$ws.Range('D12').Value = "'0.08020"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.05%  '
$ws.Range('D13').Value = '1.940.56'
$ws.Range('E13').Value = '  +0.22%  '
$ws.Range('D14').Value = "'5.350"
$ws.Range('D14').Style = 'Normal'
$ws.Range('D15').Value = "'94.54"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.26%  '
$ws.Range('D16').Value = "'14.42"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -4.44%  '
$ws.Range('D17').Value = '30.439.59'
$ws.Range('E17').Value = '  +0.31%  '
$ws.Range('D18').Value = "'254.60"
$ws.Range('D18').Style = 'Normal'
$ws.Range('D19').Value = "'0.000007941"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -3.60%  '
$ws.Range('D20').Value = "'5.767"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.05%  '
$ws.Range('D21').Value = '2.189.81'
$ws.Range('E21').Value = '  +0.05%  '
$ws.Range('D22').Value = "'1.005"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.48%  '
$ws.Range('E23').Value = '  +0.68%  '
$ws.Range('D24').Value = "'6.674"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -4.09%  '
$ws.Range('D25').Value = "'9.516"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.46%  '
$ws.Range('D26').Value = "'165.48"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.39%  '
$ws.Range('E27').Value = '  -1.37%  '
$ws.Range('D28').Value = "'0.1327"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.13%  '
$ws.Range('D29').Value = "'2.264"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -3.60%  '
$ws.Range('D30').Value = "'1.375"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.59%  '
$ws.Range('D31').Value = "'1.515"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.01%  '
$ws.Range('D32').Value = "'4.387"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.13%  '
$ws.Range('D33').Value = "'4.116"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.21%  '
$ws.Range('D34').Value = "'0.05153"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.79%  '
$ws.Range('D35').Value = "'1.280"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.99%  '
$ws.Range('D36').Value = "'0.7459"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.76%  '
$ws.Range('D37').Value = "'2.786"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.50%  '
$ws.Range('D38').Value = "'0.01950"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.05%  '
$ws.Range('D39').Value = "'2.816"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.46%  '
$ws.Range('D40').Value = "'78.74"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.88%  '
$ws.Range('D41').Value = "'6.414"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.55%  '
$ws.Range('D42').Value = "'0.4478"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.46%  '
$ws.Range('D43').Value = "'1.965"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -3.30%  '
$ws.Range('E44').Value = '  +0.51%  '
$ws.Range('D45').Value = "'0.8324"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.50%  '
$ws.Range('D46').Value = "'101.15"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.84%  '
$ws.Range('D47').Value = "'9.751"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.52%  '
$ws.Range('D48').Value = "'7.473"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.37%  '
$ws.Range('D49').Value = "'37.14"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.95%  '
$ws.Range('D50').Value = "'974.66"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +9.78%  '
$ws.Range('E51').Value = '  -0.36%  '

Write-Output "Applied cryptos update"
